$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032216756921222
$ws.Range("D2").Value = 1.04172802326125
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.04847691712512
$ws.Range("I2").Value = 1.03918694123222
$ws.Range("J2").Value = 1.037348195029395
$ws.Range("K2").Value = 1.044506717983078
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.051236653996602
$ws.Range("N2").Value = 1.016451021790868

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032985459442123
$ws.Range("D3").Value = 1.0423530425619
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.049273136818667
$ws.Range("I3").Value = 1.039373553703993
$ws.Range("J3").Value = 1.037760018603715
$ws.Range("K3").Value = 1.044942911453265
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.051844967814591
$ws.Range("N3").Value = 1.016586978700824

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033483271598796
$ws.Range("D4").Value = 1.042757773018692
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.049789130125921
$ws.Range("I4").Value = 1.039493119296797
$ws.Range("J4").Value = 1.038026217972375
$ws.Range("K4").Value = 1.04522476727446
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.052238716491283
$ws.Range("N4").Value = 1.016674853366786

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033692648246817
$ws.Range("D5").Value = 1.042927991913399
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.050006239794316
$ws.Range("I5").Value = 1.039543100079034
$ws.Range("J5").Value = 1.038138060363752
$ws.Range("K5").Value = 1.045343164556911
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.052404277575049
$ws.Range("N5").Value = 1.016711771849516

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033727809109207
$ws.Range("D6").Value = 1.042956576476938
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.050042704314679
$ws.Range("I6").Value = 1.039551475358495
$ws.Range("J6").Value = 1.038156835179421
$ws.Range("K6").Value = 1.04536303837522
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.052432077667872
$ws.Range("N6").Value = 1.016717969203615

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033486068923203
$ws.Range("D7").Value = 1.042760047218394
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.049792030427362
$ws.Range("I7").Value = 1.039493788261791
$ws.Range("J7").Value = 1.038027712683968
$ws.Range("K7").Value = 1.045226349679428
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.052240928614264
$ws.Range("N7").Value = 1.016675346768315

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032476457520828
$ws.Range("D8").Value = 1.041939187963913
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.048745839324092
$ws.Range("I8").Value = 1.039250252665192
$ws.Range("J8").Value = 1.037487429541139
$ws.Range("K8").Value = 1.044654211619611
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.051442208756227
$ws.Range("N8").Value = 1.01649698920857

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030700614233123
$ws.Range("D9").Value = 1.040495119219241
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.046908427608256
$ws.Range("I9").Value = 1.038812073161144
$ws.Range("J9").Value = 1.036533308127107
$ws.Range("K9").Value = 1.043643102730884
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.050035826544869
$ws.Range("N9").Value = 1.016181964847942

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029518990931195
$ws.Range("D10").Value = 1.039534126375553
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.04568771828459
$ws.Range("I10").Value = 1.038513927855138
$ws.Range("J10").Value = 1.035895908318901
$ws.Range("K10").Value = 1.042967145273208
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.049099057154073
$ws.Range("N10").Value = 1.015971479928228

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029007896380642
$ws.Range("D11").Value = 1.039118437542875
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.045160167451464
$ws.Range("I11").Value = 1.038383409265197
$ws.Range("J11").Value = 1.035619610200851
$ws.Range("K11").Value = 1.042674018899217
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.0486936416431
$ws.Range("N11").Value = 1.015880231921481

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028818138465566
$ws.Range("D12").Value = 1.038964098341674
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.044964367591444
$ws.Range("I12").Value = 1.038334716408619
$ws.Range("J12").Value = 1.035516937008945
$ws.Range("K12").Value = 1.042565075182442
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.048543085850918
$ws.Range("N12").Value = 1.015846322764163

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028858838290924
$ws.Range("D13").Value = 1.038997201600401
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.045006360232507
$ws.Range("I13").Value = 1.038345170791338
$ws.Range("J13").Value = 1.035538962717359
$ws.Range("K13").Value = 1.042588446831545
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.048575379031683
$ws.Range("N13").Value = 1.015853597091739

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028992209180486
$ws.Range("D14").Value = 1.039105678456003
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.045143979379235
$ws.Range("I14").Value = 1.038379388627276
$ws.Range("J14").Value = 1.035611124088297
$ws.Range("K14").Value = 1.042665014865374
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.048681195959187
$ws.Range("N14").Value = 1.015877429295355

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029074394752004
$ws.Range("D15").Value = 1.039172523448632
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.045228791813259
$ws.Range("I15").Value = 1.038400443240382
$ws.Range("J15").Value = 1.035655579328804
$ws.Range("K15").Value = 1.042712182582272
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.048746397762761
$ws.Range("N15").Value = 1.01589211105051

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029552922464434
$ws.Range("D16").Value = 1.039561723443904
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.045722751858204
$ws.Range("I16").Value = 1.038522560105412
$ws.Range("J16").Value = 1.035914239107319
$ws.Range("K16").Value = 1.042986590118016
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.049125967859071
$ws.Range("N16").Value = 1.015977533548364

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029853240423146
$ws.Range("D17").Value = 1.039805974109423
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.046032875684806
$ws.Range("I17").Value = 1.038598781295807
$ws.Range("J17").Value = 1.036076410187602
$ws.Range("K17").Value = 1.043158603972701
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.049364120102594
$ws.Range("N17").Value = 1.016031088581152

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030028464345627
$ws.Range("D18").Value = 1.039948482500179
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.046213864360469
$ws.Range("I18").Value = 1.038643102788402
$ws.Range("J18").Value = 1.036170972823939
$ws.Range("K18").Value = 1.043258894934585
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.049503050596626
$ws.Range("N18").Value = 1.016062316006443

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030088220201107
$ws.Range("D19").Value = 1.039997081092975
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.046275593497759
$ws.Range("I19").Value = 1.038658192001982
$ws.Range("J19").Value = 1.036203211279704
$ws.Range("K19").Value = 1.043293084425921
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.049550425697884
$ws.Range("N19").Value = 1.016072961987823

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029821013597635
$ws.Range("D20").Value = 1.039779764058482
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.045999592109534
$ws.Range("I20").Value = 1.038590617651771
$ws.Range("J20").Value = 1.03605901375703
$ws.Range("K20").Value = 1.043140152828277
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.049338566506286
$ws.Range("N20").Value = 1.016025343697183

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028952932406918
$ws.Range("D21").Value = 1.039073732884513
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.045103449647057
$ws.Range("I21").Value = 1.038369318185002
$ws.Range("J21").Value = 1.035589875561152
$ws.Range("K21").Value = 1.042642469230139
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.048650034570395
$ws.Range("N21").Value = 1.015870411732097

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02840763005354
$ws.Range("D22").Value = 1.038630206677987
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.044540912973241
$ws.Range("I22").Value = 1.038228949731917
$ws.Range("J22").Value = 1.035294657153773
$ws.Range("K22").Value = 1.042329188936922
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.048217322364259
$ws.Range("N22").Value = 1.015772909903769

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028696657649901
$ws.Range("D23").Value = 1.038865291307537
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.044839037947168
$ws.Range("I23").Value = 1.038303477885723
$ws.Range("J23").Value = 1.035451181536786
$ws.Range("K23").Value = 1.042495299047374
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.04844669220804
$ws.Range("N23").Value = 1.015824605848244

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029835575347697
$ws.Range("D24").Value = 1.039791607125531
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.046014631219848
$ws.Range("I24").Value = 1.03859430687464
$ws.Range("J24").Value = 1.036066874544985
$ws.Range("K24").Value = 1.043148490235586
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.049350113013957
$ws.Range("N24").Value = 1.016027939594903

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031159319599789
$ws.Range("D25").Value = 1.040868150736175
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.047382705712767
$ws.Range("I25").Value = 1.038926418893279
$ws.Range("J25").Value = 1.03678020934378
$ws.Range("K25").Value = 1.043904836985776
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.050399273295888
$ws.Range("N25").Value = 1.016263490695101
